$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 206, shifting existing rows 206:284 down to 207:285
$ws.Rows.Item(206).EntireRow.Insert()

# Populate the newly inserted row 206 with the new record's data
$ws.Range("A206").Value = 3
$ws.Range("B206").Value = "Femacal de La Calera"
$ws.Range("C206").Value = "Coquimbo"
$ws.Range("D206").Value = 45146
$ws.Range("E206").Value = 5
$ws.Range("F206").Value = 100112026
$ws.Range("G206").Value = "Haba"
$ws.Range("H206").Value = "Sin especificar"
$ws.Range("I206").Value = "Primera"
$ws.Range("J206").Value = 40
$ws.Range("K206").Value = 15000
$ws.Range("L206").Value = 15000
$ws.Range("M206").Value = 15000
$ws.Range("N206").Value = '$/saco 25 kilos'
$ws.Range("O206").Value = "Provincia de Limarí"
$ws.Range("P206").Value = 600
$ws.Range("Q206").Value = 25
$ws.Range("R206").Value = "Hortaliza"
